$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended below the existing table (rows 23-26).
$newRows = @(
    @{ Row = 23; Ref = "3M PBO 1.5 (FDP Branch)"; Nbr = 12; Taille = 12 },
    @{ Row = 24; Ref = "3M PBO T1 (FDP Branch)";  Nbr = 12; Taille = 12 },
    @{ Row = 25; Ref = "3M PBO T0 (Branch)";       Nbr = 4;  Taille = 12 },
    @{ Row = 26; Ref = "3M BPEO T2 (CDP)";         Nbr = 14; Taille = 24 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Ref
    $ws.Cells.Item($row, 2).Value = $r.Nbr
    $ws.Cells.Item($row, 3).Value = $r.Taille
}

# The new rows use an integer number format, which introduces a new cell
# style (numFmtId 1 "0") in the workbook's style table.
$ws.Range("A23:C26").NumberFormat = "0"

# Update selection / scroll to match the author's final view state.
$ws.Range("A23:C26").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
